$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update input values
$ws.Range("B1").Value = 30000000
$ws.Range("F1").Value = 10
$ws.Range("F2").Value = 10
$ws.Range("B3").Value = 480

# B2, B4, B5 previously held formulas; they are now plain values
$ws.Range("B2").Value = 100
$ws.Range("B4").Value = 1000
$ws.Range("B5").Value = 1000

# Update the active selection to match the new state
$ws.Range("E8").Select()
